$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8bd46503e3696386edea40032b2ad9c3bff9ef03/e2e/6bb944c4-1810-40d7-989b-afb5716f321b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/74c2d705d5af1a0f53935fa7299a8c42141dfa72/e2e/6bb944c4-1810-40d7-989b-afb5716f321b.md."

# ---------- zh-cn sheet ----------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I8").Value = "6bb944c4-1810-40d7-989b-afb5716f321b.md"
$wsZh.Range("J8").Value = "6bb944c4-1810-40d7-989b-afb5716f321b.6e2519e4385f6cf436f99306e52624509aac9092.zh-cn.xlf"
$wsZh.Range("K8").Value = "2016-08-22 02:53:16"
$wsZh.Range("P8").Value = $errorDetail

$wsZh.Range("I8").Font.Underline = $true
$wsZh.Range("I8").Font.Color = 15570276

$wsZh.Hyperlinks.Add($wsZh.Range("I8"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/8bd46503e3696386edea40032b2ad9c3bff9ef03/e2e/6bb944c4-1810-40d7-989b-afb5716f321b.md", [Type]::Missing, [Type]::Missing, "6bb944c4-1810-40d7-989b-afb5716f321b.md")

$wsZh.Range("I8").Value = "6bb944c4-1810-40d7-989b-afb5716f321b.md"
$wsZh.Range("I8").Font.Underline = $true
$wsZh.Range("I8").Font.Color = 15570276

$wsZh.Columns.Item(16).ColumnWidth = 40

# ---------- de-de sheet ----------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I8").Value = "6bb944c4-1810-40d7-989b-afb5716f321b.md"
$wsDe.Range("J8").Value = "6bb944c4-1810-40d7-989b-afb5716f321b.6e2519e4385f6cf436f99306e52624509aac9092.de-de.xlf"
$wsDe.Range("K8").Value = "2016-08-22 02:53:22"
$wsDe.Range("P8").Value = $errorDetail

$wsDe.Hyperlinks.Add($wsDe.Range("I8"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8bd46503e3696386edea40032b2ad9c3bff9ef03/e2e/6bb944c4-1810-40d7-989b-afb5716f321b.md", [Type]::Missing, [Type]::Missing, "6bb944c4-1810-40d7-989b-afb5716f321b.md")

$wsDe.Range("I8").Value = "6bb944c4-1810-40d7-989b-afb5716f321b.md"
$wsDe.Range("I8").Font.Underline = $true
$wsDe.Range("I8").Font.Color = 15570276

$wsDe.Columns.Item(16).ColumnWidth = 40

Write-Host "Applied handback report changes"
